$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D data cells remain text (avoid Excel auto-converting numeric-looking strings)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.444.45'
$ws.Range('E2').Value = '  +3.18%  '
$ws.Range('D3').Value = '1.605.51'
$ws.Range('E3').Value = '  +2.83%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '212.25'
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('E6').Value = '  +6.98%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').Value = '26.71'
$ws.Range('E8').Value = '  +7.58%  '
$ws.Range('D9').Value = '43.43'
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('E10').Value = '  +2.88%  '
$ws.Range('D11').Value = '0.0600'
$ws.Range('E11').Value = '  +2.63%  '
$ws.Range('E12').Value = '  +1.65%  '
$ws.Range('D13').Value = '1.835.27'
$ws.Range('E13').Value = '  +2.72%  '
$ws.Range('D14').Value = '1.605.94'
$ws.Range('E14').Value = '  +2.82%  '
$ws.Range('D15').Value = '29.447.34'
$ws.Range('E15').Value = '  +3.11%  '
$ws.Range('D16').Value = '0.534'
$ws.Range('E16').Value = '  +4.14%  '
$ws.Range('D17').Value = '3.71'
$ws.Range('E17').Value = '  +2.29%  '
$ws.Range('D18').Value = '63.08'
$ws.Range('E18').Value = '  +3.13%  '
$ws.Range('D19').Value = '241.60'
$ws.Range('E19').Value = '  +5.14%  '
$ws.Range('D20').Value = '7.65'
$ws.Range('E20').Value = '  +4.17%  '
$ws.Range('D21').Value = '0.0₃0688'
$ws.Range('E21').Value = '  +2.04%  '
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('E23').Value = '  +2.29%  '
$ws.Range('D24').Value = '9.16'
$ws.Range('E24').Value = '  +2.29%  '
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('D26').Value = '154.64'
$ws.Range('E26').Value = '  +2.98%  '
$ws.Range('E27').Value = '  +5.26%  '
$ws.Range('D28').Value = '15.28'
$ws.Range('E28').Value = '  +3.46%  '
$ws.Range('E29').Value = '  +2.69%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('E31').Value = '  +2.60%  '
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('E33').Value = '  +1.73%  '
$ws.Range('D34').Value = '3.10'
$ws.Range('E34').Value = '  +4.49%  '
$ws.Range('D35').Value = '1.412.68'
$ws.Range('E35').Value = '  +1.79%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E37').Value = '  +3.43%  '
$ws.Range('D38').Value = '2.82'
$ws.Range('E38').Value = '  +4.71%  '
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('E40').Value = '  +2.70%  '
$ws.Range('E41').Value = '  +3.83%  '
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('D43').Value = '0.0489'
$ws.Range('E43').Value = '  +6.09%  '
$ws.Range('B44').Value = 'BitcoinSV'
$ws.Range('C44').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D44').Value = '52.87'
$ws.Range('E44').Value = '  +23.07%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').Value = '0.796'
$ws.Range('E45').Value = '  +3.22%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value = '0.998'
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('D47').Value = '65.61'
$ws.Range('E47').Value = '  +2.95%  '
$ws.Range('D48').Value = '1.745.52'
$ws.Range('E48').Value = '  +2.70%  '
$ws.Range('E49').Value = '  +0.54%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').Value = '0.854'
$ws.Range('E50').Value = '  -1.82%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = '86.64'
$ws.Range('E51').Value = '  +1.81%  '
